$wb = $excel.ActiveWorkbook

# --- Update "Model Profile" sheet (F4:G15) ---
$wsProfile = $wb.Worksheets.Item("Model Profile")

$profileUpdates = @(
    @{Row=4;  F=15.76137065887451;  G=230.6125734453858},
    @{Row=5;  F=21.6904091835022;   G=273.1469391460485},
    @{Row=6;  F=24.83202219009399;  G=381.4970774020932},
    @{Row=7;  F=26.2183928489685;   G=429.5303772916355},
    @{Row=8;  F=25.70167064666748;  G=90.43009740517203},
    @{Row=9;  F=29.12822961807251;  G=137.2221240733107},
    @{Row=10; F=34.91165161132812;  G=182.5665340188121},
    @{Row=11; F=52.65275239944457;  G=246.5175600174442},
    @{Row=12; F=22.60903596878051;  G=102.7385153817363},
    @{Row=13; F=25.35673141479492;  G=179.7882116347376},
    @{Row=14; F=29.36871290206909;  G=259.6441140147517},
    @{Row=15; F=39.58753108978271;  G=342.2610144325105}
)

foreach ($u in $profileUpdates) {
    $wsProfile.Cells.Item($u.Row, 6).Value = $u.F
    $wsProfile.Cells.Item($u.Row, 7).Value = $u.G
}

# --- Update "Model Raw Profile" sheet (C2 JSON string) ---
$wsRaw = $wb.Worksheets.Item("Model Raw Profile")

$newJson = '{"Sentimental-bert24-2/tokenizer": {"CPU1": {"THROUGHPUT": [[1, 230.61257344538583], [2, 273.1469391460485], [4, 381.49707740209317], [8, 429.53037729163555]], "LATENCY": [[1, 15.761370658874512], [2, 21.690409183502197], [4, 24.832022190093994], [8, 26.218392848968495]]}}, "Sentimental-bert24-2/bert24_p2_stage0": {"Tesla P40": {"THROUGHPUT": [[1, 90.43009740517203], [2, 137.22212407331068], [4, 182.56653401881212], [8, 246.51756001744423]], "LATENCY": [[1, 25.70167064666748], [2, 29.128229618072506], [4, 34.91165161132812], [8, 52.652752399444566]]}}, "Sentimental-bert24-2/bert24-p2-stage1": {"Tesla P40": {"THROUGHPUT": [[1, 102.7385153817363], [2, 179.78821163473765], [4, 259.6441140147517], [8, 342.26101443251054]], "LATENCY": [[1, 22.609035968780514], [2, 25.356731414794922], [4, 29.36871290206909], [8, 39.587531089782715]]}}}'

$wsRaw.Range("C2").Value = $newJson
